# [Kadastro App] Yeni kayit eklendi: 2532
# Adds the new record row to both the master "Kayitlar" log sheet and the
# per-district "Erdemli" sheet (one row appended at the bottom of each table).

$wb = $excel.ActiveWorkbook

$recordNo   = "2532"
$date       = "2025-08-07"
$unit       = "Erdemli"
$parcelCnt  = "1"
$job        = "MAKS"
$personnel  = "ÖZKAN AKBAŞ (Mühendis), CEMAL TİMUROĞLU (K.Teknisyeni)"

# --- "Kayitlar" master sheet: append after the last used row ---------------
$wsKayitlar = $wb.Worksheets.Item("Kayitlar")
$lastRowKayitlar = $wsKayitlar.Cells.Item($wsKayitlar.Rows.Count, 1).End(-4162).Row
$newRowKayitlar = $lastRowKayitlar + 1

$wsKayitlar.Cells.Item($newRowKayitlar, 1).NumberFormat = "@"
$wsKayitlar.Cells.Item($newRowKayitlar, 1).Value = "'" + $recordNo
$wsKayitlar.Cells.Item($newRowKayitlar, 2).NumberFormat = "@"
$wsKayitlar.Cells.Item($newRowKayitlar, 2).Value = "'" + $date
$wsKayitlar.Cells.Item($newRowKayitlar, 3).NumberFormat = "@"
$wsKayitlar.Cells.Item($newRowKayitlar, 3).Value = "'" + $unit
$wsKayitlar.Cells.Item($newRowKayitlar, 4).NumberFormat = "@"
$wsKayitlar.Cells.Item($newRowKayitlar, 4).Value = "'" + $parcelCnt
$wsKayitlar.Cells.Item($newRowKayitlar, 5).NumberFormat = "@"
$wsKayitlar.Cells.Item($newRowKayitlar, 5).Value = "'" + $job
$wsKayitlar.Cells.Item($newRowKayitlar, 6).NumberFormat = "@"
$wsKayitlar.Cells.Item($newRowKayitlar, 6).Value = "'" + $personnel

# --- "Erdemli" district sheet: append after the last used row --------------
$wsErdemli = $wb.Worksheets.Item("Erdemli")
$lastRowErdemli = $wsErdemli.Cells.Item($wsErdemli.Rows.Count, 1).End(-4162).Row
$newRowErdemli = $lastRowErdemli + 1

$wsErdemli.Cells.Item($newRowErdemli, 1).NumberFormat = "@"
$wsErdemli.Cells.Item($newRowErdemli, 1).Value = "'" + $recordNo
$wsErdemli.Cells.Item($newRowErdemli, 2).NumberFormat = "@"
$wsErdemli.Cells.Item($newRowErdemli, 2).Value = "'" + $date
$wsErdemli.Cells.Item($newRowErdemli, 3).NumberFormat = "@"
$wsErdemli.Cells.Item($newRowErdemli, 3).Value = "'" + $unit
$wsErdemli.Cells.Item($newRowErdemli, 4).NumberFormat = "@"
$wsErdemli.Cells.Item($newRowErdemli, 4).Value = "'" + $parcelCnt
$wsErdemli.Cells.Item($newRowErdemli, 5).NumberFormat = "@"
$wsErdemli.Cells.Item($newRowErdemli, 5).Value = "'" + $job
$wsErdemli.Cells.Item($newRowErdemli, 6).NumberFormat = "@"
$wsErdemli.Cells.Item($newRowErdemli, 6).Value = "'" + $personnel
